$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for data rows 2-33 is updated from 2023-10-13 (45212)
# to 2023-10-22 (45221) for every row in the sheet.
$ws.Range("C2:C33").Value = (Get-Date -Year 2023 -Month 10 -Day 22 -Hour 0 -Minute 0 -Second 0).Date
